$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the price cells whose new values look numeric as Text,
# so Excel stores them as strings (matching the source data) instead of
# silently converting them to numbers / dropping significant trailing zeros.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated coin data (rows 2-51 of the table).
$ws.Range("D2").Value = '30.345.75'
$ws.Range("E2").Value = '  +1.16%  '
$ws.Range("D3").Value = '1.923.14'
$ws.Range("E3").Value = '  +0.71%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '0.8093'
$ws.Range("E5").Value = '  +2.29%  '
$ws.Range("D6").Value = '244.52'
$ws.Range("E6").Value = '  +1.07%  '
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = '0.3260'
$ws.Range("E8").Value = '  +3.04%  '
$ws.Range("D9").Value = '27.16'
$ws.Range("E9").Value = '  +3.31%  '
$ws.Range("D10").Value = '0.07266'
$ws.Range("E10").Value = '  +5.47%  '
$ws.Range("D11").Value = '0.7931'
$ws.Range("D12").Value = '0.08106'
$ws.Range("E12").Value = '  +1.41%  '
$ws.Range("D13").Value = '1.941.33'
$ws.Range("E13").Value = '  +1.93%  '
$ws.Range("D14").Value = '5.419'
$ws.Range("E14").Value = '  +4.45%  '
$ws.Range("D15").Value = '94.38'
$ws.Range("E15").Value = '  +1.48%  '
$ws.Range("D16").Value = '30.378.30'
$ws.Range("E16").Value = '  +1.27%  '
$ws.Range("D17").Value = '14.34'
$ws.Range("E17").Value = '  +2.95%  '
$ws.Range("D18").Value = '6.102'
$ws.Range("E18").Value = '  +4.05%  '
$ws.Range("D19").Value = '251.38'
$ws.Range("E19").Value = '  +2.35%  '
$ws.Range("D20").Value = '0.000007865'
$ws.Range("E20").Value = '  +1.72%  '
$ws.Range("D21").Value = '2.181.88'
$ws.Range("E21").Value = '  +1.11%  '
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("D23").Value = '8.015'
$ws.Range("E23").Value = '  +17.27%  '
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").Value = '0.1671'
$ws.Range("E25").Value = '  +19.56%  '
$ws.Range("D26").Value = '9.526'
$ws.Range("E26").Value = '  +3.22%  '
$ws.Range("D27").Value = '167.68'
$ws.Range("E27").Value = '  -0.22%  '
$ws.Range("D28").Value = '19.12'
$ws.Range("E28").Value = '  +1.21%  '
$ws.Range("D29").Value = '2.159'
$ws.Range("E29").Value = '  +6.39%  '
$ws.Range("D30").Value = '1.376'
$ws.Range("E30").Value = '  +0.82%  '
$ws.Range("E31").Value = '  +2.19%  '
$ws.Range("D32").Value = '4.352'
$ws.Range("E32").Value = '  +0.87%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.05649'
$ws.Range("E33").Value = '  +2.03%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '4.140'
$ws.Range("E34").Value = '  +1.58%  '
$ws.Range("D35").Value = '1.303'
$ws.Range("E35").Value = '  +3.91%  '
$ws.Range("D36").Value = '0.7463'
$ws.Range("E36").Value = '  +1.71%  '
$ws.Range("E37").Value = '  +0.14%  '
$ws.Range("D38").Value = '2.727'
$ws.Range("E38").Value = '  +0.21%  '
$ws.Range("D39").Value = '0.01961'
$ws.Range("E39").Value = '  +2.03%  '
$ws.Range("D40").Value = '2.823'
$ws.Range("E40").Value = '  +1.43%  '
$ws.Range("D41").Value = '0.4506'
$ws.Range("E41").Value = '  +2.10%  '
$ws.Range("D42").Value = '74.29'
$ws.Range("E42").Value = '  +2.75%  '
$ws.Range("D43").Value = '5.994'
$ws.Range("E43").Value = '  -2.45%  '
$ws.Range("D44").Value = '0.8571'
$ws.Range("D45").Value = '1.930'
$ws.Range("E45").Value = '  +2.71%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '1.036.50'
$ws.Range("E47").Value = '  +4.97%  '
$ws.Range("D48").Value = '103.18'
$ws.Range("E48").Value = '  +2.63%  '
$ws.Range("D49").Value = '9.964'
$ws.Range("E49").Value = '  +2.14%  '
$ws.Range("D50").Value = '7.654'
$ws.Range("E50").Value = '  +1.41%  '
$ws.Range("B51").Value = 'SynthetixNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D51").Value = '3.039'
$ws.Range("E51").Value = '  +8.86%  '
